# Move the "PDE-CL" title-graphic slide (currently the 7th slide) to the
# very front of the deck, ahead of the original "Partial Differential
# Equation Constrained Layer (PDE-CL) for Inverse Problems and Design"
# title slide. Everything else keeps its relative order.
$p = $ppt.ActivePresentation
$p.Slides.Item(7).MoveTo(1)
